$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four data rows (2-5) are rotated: new row2 = old row4, new row3 = old row5,
# new row4 = old row2, new row5 = old row3. Columns A,B,C,E,F,G,H,O,R stay the same
# for each row; columns D,I,J,K,L,M,N,P,Q carry the rotated values.

$ws.Range("D2").Value = 44267
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = 1650
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("P2").Value = 275
$ws.Range("Q2").Value = 6

$ws.Range("D3").Value = 44623
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1900
$ws.Range("N3").Value = "$/paquete"
$ws.Range("P3").Value = 1900
$ws.Range("Q3").Value = 1

$ws.Range("D4").Value = 44377
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2800
$ws.Range("M4").Value = 2364
$ws.Range("N4").Value = "$/docena de matas"
$ws.Range("P4").Value = 394
$ws.Range("Q4").Value = 6

$ws.Range("D5").Value = 44370
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1200
$ws.Range("M5").Value = 1080
$ws.Range("N5").Value = "$/docena de matas"
$ws.Range("P5").Value = 180
$ws.Range("Q5").Value = 6
